$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates mirroring the refreshed coinranking.com snapshot.
# Values in column D that are plain decimal numbers must be force-written
# as text (matching the original inlineStr/shared-string cell type) so Excel
# does not silently convert them to floating point numbers.

$ws.Range("D2").Value = "68.010.10"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.421.49"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.155"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.24%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.325"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "67.941.07"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000168"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "334.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "0.0₃0800"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "419.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.102"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.291"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "127.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0713"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.472"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.554"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.07%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0423"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
